# Auto-generated Excel COM-interop script to apply cryptos.xlsx price/volume updates
# (Sun Dec 10 17:30:14 UTC 2023 GitHub Actions cryptos-list refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '43.870.23'
$ws.Range('E2').Value = '  -0.22%  '
$ws.Range('D3').Value = '2.346.77'
$ws.Range('E3').Value = '  -0.11%  '
$ws.Range('E4').Value = '  +0.04%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.22'
$ws.Range('D5').ClearFormats()
$ws.Range('E5').Value = '  +0.25%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.665'
$ws.Range('D6').ClearFormats()
$ws.Range('E6').Value = '  -1.58%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '72.97'
$ws.Range('D7').ClearFormats()
$ws.Range('E7').Value = '  -1.16%  '
$ws.Range('E8').Value = '  +0.04%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.595'
$ws.Range('D9').ClearFormats()
$ws.Range('E9').Value = '  +0.54%  '
$ws.Range('B10').Value = 'Dogecoin'
$ws.Range('C10').Value = 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge'
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.101'
$ws.Range('D10').ClearFormats()
$ws.Range('E10').Value = '  +0.76%  '
$ws.Range('B11').Value = 'OKB'
$ws.Range('C11').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '60.95'
$ws.Range('D11').ClearFormats()
$ws.Range('E11').Value = '  +6.61%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '33.42'
$ws.Range('D12').ClearFormats()
$ws.Range('E12').Value = '  +4.75%  '
$ws.Range('E13').Value = '  +0.73%  '
$ws.Range('E14').Value = '  +0.96%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '16.10'
$ws.Range('D15').ClearFormats()
$ws.Range('E15').Value = '  -2.18%  '
$ws.Range('E16').Value = '  +0.81%  '
$ws.Range('D17').Value = '2.332.21'
$ws.Range('E17').Value = '  -0.78%  '
$ws.Range('D18').Value = '43.733.67'
$ws.Range('E18').Value = '  -0.39%  '
$ws.Range('E19').Value = '  +0.14%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '77.82'
$ws.Range('D20').ClearFormats()
$ws.Range('E20').Value = '  +1.60%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '6.56'
$ws.Range('D21').ClearFormats()
$ws.Range('E21').Value = '  -1.11%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '251.76'
$ws.Range('D22').ClearFormats()
$ws.Range('E22').Value = '  -1.56%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '3.82'
$ws.Range('D23').ClearFormats()
$ws.Range('E23').Value = '  +3.26%  '
$ws.Range('E25').Value = '  -3.22%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '2.49'
$ws.Range('D26').ClearFormats()
$ws.Range('E26').Value = '  +0.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '10.40'
$ws.Range('D27').ClearFormats()
$ws.Range('E27').Value = '  -2.10%  '
$ws.Range('E28').Value = '  +0.81%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '175.65'
$ws.Range('D29').ClearFormats()
$ws.Range('E29').Value = '  +0.51%  '
$ws.Range('E30').Value = '  -2.23%  '
$ws.Range('E31').Value = '  +0.49%  '
$ws.Range('E32').Value = '  -2.05%  '
$ws.Range('E33').Value = '  -2.14%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '5.03'
$ws.Range('D34').ClearFormats()
$ws.Range('E34').Value = '  -4.15%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.33'
$ws.Range('D35').ClearFormats()
$ws.Range('E35').Value = '  -0.21%  '
$ws.Range('E36').Value = '  +0.99%  '
$ws.Range('B37').Value = 'LidoDAOToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo'
$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '2.38'
$ws.Range('D37').ClearFormats()
$ws.Range('E37').Value = '  +1.47%  '
$ws.Range('B38').Value = 'THORChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/ybmU-kKU+thorchain-rune'
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '6.40'
$ws.Range('D38').ClearFormats()
$ws.Range('E38').Value = '  +1.42%  '
$ws.Range('E39').Value = '  -3.39%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '5.39'
$ws.Range('D40').ClearFormats()
$ws.Range('E40').Value = '  +15.40%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '64.84'
$ws.Range('D41').ClearFormats()
$ws.Range('E41').Value = '  +12.93%  '
$ws.Range('E42').Value = '  +2.04%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '9.08'
$ws.Range('D43').ClearFormats()
$ws.Range('E43').Value = '  +0.67%  '
$ws.Range('E45').Value = '  -1.51%  '
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('E47').Value = '  -1.44%  '
$ws.Range('B48').Value = 'TerraClassic'
$ws.Range('C48').Value = 'https://coinranking.com/coin/AaQUAs2Mc+terraclassic-lunc'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '0.000216'
$ws.Range('D48').ClearFormats()
$ws.Range('E48').Value = '  +18.53%  '
$ws.Range('B49').Value = 'NEARProtocol'
$ws.Range('C49').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.42'
$ws.Range('D49').ClearFormats()
$ws.Range('E49').Value = '  -1.87%  '
$ws.Range('B50').Value = 'ARBITRUM'
$ws.Range('C50').Value = 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb'
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '1.15'
$ws.Range('D50').ClearFormats()
$ws.Range('E50').Value = '  -1.70%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '97.42'
$ws.Range('D51').ClearFormats()
$ws.Range('E51').Value = '  -2.25%  '
